# 10Th - MB for single stock and added new group
#
# This report tracks a rolling window of "market beat rank" snapshots per
# analyst/firm (column headers are dates, most-recent first in column B).
# A new snapshot date (Jun_27) is being added, which pushes the existing
# date columns one slot to the right. The underlying 6/26/2018 downgrade
# event (previously only implied) is also recorded explicitly for
# BidaskClub, and two new research firms are appended as new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new (most-recent) date columns in front of the existing
#     data. This shifts the old B:E (Jun_17, Jun_15, Jun_13, Jun_10)
#     columns right to E:H, values/styles intact.
$ws.Range("B:D").Insert()

# --- New header row values for the inserted date columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- Every analyst/firm row (2-27) gets "UN" (unchanged) in the three
#     newly inserted columns by default.
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("B" + $r).Value = "UN"
    $ws.Range("C" + $r).Value = "UN"
    $ws.Range("D" + $r).Value = "UN"
}

# --- BidaskClub (row 22) had a rating change on 6/26/2018; record the
#     detail in the new columns and highlight the triggering cell.
$downgradeNote = "6/26/2018,Downgrades,Strong-Buy -> Buy,"
$ws.Range("B22").Value = $downgradeNote
$ws.Range("C22").Value = $downgradeNote
$ws.Range("D22").Value = $downgradeNote
$ws.Range("B22").Interior.ColorIndex = 45

# --- Append two newly-tracked research firms as new rows.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
